$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 216.33333
$ws.Range("I9").Value = 170
$ws.Range("K9").Value = 170
$ws.Range("M9").Value = -1
# Row 40
$ws.Range("H40").Value = 3968.8
$ws.Range("J40").Value = 3968.8
$ws.Range("L40").Value = 3968.8
$ws.Range("N40").Value = -4318.8
# Row 53
$ws.Range("H53").Value = 309.14816
$ws.Range("J53").Value = 153.8
$ws.Range("L53").Value = 153.8
$ws.Range("N53").Value = -1427.8
# Row 138
$ws.Range("H138").Value = 3189.6
$ws.Range("J138").Value = 3900
$ws.Range("L138").Value = 11700
$ws.Range("N138").Value = -21980

# ---- ARM sheet ----
$ws = $wb.Worksheets.Item("ARM")
# Row 24
$ws.Range("H24").Value = 42392.332
$ws.Range("J24").Value = 42392.332
$ws.Range("L24").Value = 42392.332
$ws.Range("N24").Value = -43140.332
# Row 63
$ws.Range("H63").Value = 3313.25
$ws.Range("I63").Value = 2451.8
$ws.Range("K63").Value = 2451.8
$ws.Range("M63").Value = -1765.8
# Row 66
$ws.Range("H66").Value = 3313.25
$ws.Range("I66").Value = 2451.8
$ws.Range("K66").Value = 12259
$ws.Range("M66").Value = -8827
# Row 100
$ws.Range("H100").Value = 42392.332
$ws.Range("J100").Value = 42392.332
$ws.Range("L100").Value = 42392.332
$ws.Range("N100").Value = -44556.332
# Row 104
$ws.Range("H104").Value = 23408
$ws.Range("J104").Value = 23408
$ws.Range("L104").Value = 23408
$ws.Range("N104").Value = -30396
# Row 106
$ws.Range("H106").Value = 49456.332
$ws.Range("J106").Value = 49456.332
$ws.Range("L106").Value = 49456.332
$ws.Range("N106").Value = -51980.332
# Row 110
$ws.Range("H110").Value = 2115.875
$ws.Range("I110").Value = 1032.1818
$ws.Range("K110").Value = 1032.1818
$ws.Range("M110").Value = 1012.8182
# Row 132
$ws.Range("H132").Value = 3821.5557
$ws.Range("I132").Value = 3811.75
$ws.Range("K132").Value = 11435.25
$ws.Range("M132").Value = -8905.25

# ---- BSM sheet ----
$ws = $wb.Worksheets.Item("BSM")
# Row 64
$ws.Range("H64").Value = 1082.125
$ws.Range("I64").Value = 1354.5
$ws.Range("J64").Value = 918.7
$ws.Range("K64").Value = 1354.5
$ws.Range("L64").Value = 918.7
$ws.Range("M64").Value = -1129.5
$ws.Range("N64").Value = -1368.7
# Row 67
$ws.Range("H67").Value = 1082.125
$ws.Range("I67").Value = 1354.5
$ws.Range("J67").Value = 918.7
$ws.Range("K67").Value = 1354.5
$ws.Range("L67").Value = 918.7
$ws.Range("M67").Value = -574.5
$ws.Range("N67").Value = -2478.7
# Row 100
$ws.Range("H100").Value = 39829.5
$ws.Range("J100").Value = 39829.5
$ws.Range("L100").Value = 39829.5
$ws.Range("N100").Value = -41993.5

# ---- CRP sheet ----
$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 90
$ws.Range("J4").Value = 90
$ws.Range("L4").Value = 90
$ws.Range("N4").Value = -314
# Row 64
$ws.Range("H64").Value = 30000
$ws.Range("J64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30496
# Row 67
$ws.Range("H67").Value = 30000
$ws.Range("J67").Value = 30000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31716
# Row 96
$ws.Range("H96").Value = 20869.4
$ws.Range("J96").Value = 20869.4
$ws.Range("L96").Value = 20869.4
$ws.Range("N96").Value = -26361.4
# Row 99
$ws.Range("H99").Value = 4713.857
$ws.Range("J99").Value = 5624.5
$ws.Range("L99").Value = 5624.5
$ws.Range("N99").Value = -8620.5
# Row 126
$ws.Range("H126").Value = 4713.857
$ws.Range("J126").Value = 5624.5
$ws.Range("L126").Value = 16873.5
$ws.Range("N126").Value = -21813.5
# Row 132
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

# ---- CUL sheet ----
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 251.8
$ws.Range("J12").Value = 327
$ws.Range("L12").Value = 981
$ws.Range("N12").Value = -1327
# Row 97
$ws.Range("H97").Value = 779.7143
$ws.Range("I97").Value = 1183.3334
$ws.Range("K97").Value = 3550.0002
$ws.Range("M97").Value = -3054.0002
# Row 102
$ws.Range("H102").Value = 2599.5
$ws.Range("I102").Value = 2599.5
$ws.Range("K102").Value = 7798.5
$ws.Range("M102").Value = -5364.5
# Row 119
$ws.Range("H119").Value = 770
$ws.Range("I119").Value = 770
$ws.Range("K119").Value = 2310
$ws.Range("M119").Value = 2528
# Row 121
$ws.Range("H121").Value = 2723
$ws.Range("I121").Value = 2000
$ws.Range("K121").Value = 6000
$ws.Range("M121").Value = -4690

# ---- GSM sheet ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 8671.308000000001
$ws.Range("I80").Value = 8249.714
$ws.Range("K80").Value = 8249.714
$ws.Range("M80").Value = -7251.714
# Row 83
$ws.Range("H83").Value = 8671.308000000001
$ws.Range("I83").Value = 8249.714
$ws.Range("K83").Value = 41248.57
$ws.Range("M83").Value = -36256.57

# ---- LTW sheet ----
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 2746.4443
$ws.Range("I46").Value = 1990
$ws.Range("J46").Value = 3124.6667
$ws.Range("K46").Value = 1990
$ws.Range("L46").Value = 3124.6667
$ws.Range("N46").Value = -3500.6667
$ws.Range("M46").Value = -1802
# Row 95
$ws.Range("H95").Value = 17125
$ws.Range("J95").Value = 17125
$ws.Range("L95").Value = 17125
$ws.Range("N95").Value = -22617

# ---- WVR sheet ----
$ws = $wb.Worksheets.Item("WVR")
# Row 68
$ws.Range("H68").Value = 70000
$ws.Range("J68").Value = 70000
$ws.Range("L68").Value = 70000
$ws.Range("N68").Value = -71622
# Row 71
$ws.Range("H71").Value = 70000
$ws.Range("J71").Value = 70000
$ws.Range("L71").Value = 210000
$ws.Range("N71").Value = -218112
# Row 81
$ws.Range("H81").Value = 698
$ws.Range("I81").Value = 698
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 1396
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -335
$ws.Range("N81").ClearContents()
# Row 84
$ws.Range("H84").Value = 698
$ws.Range("I84").Value = 698
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 6980
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -1676
$ws.Range("N84").ClearContents()
# Row 104
$ws.Range("H104").Value = 29499.666
$ws.Range("J104").Value = 29499.666
$ws.Range("L104").Value = 29499.666
$ws.Range("N104").Value = -36487.666

Write-Host "All updates applied"
